$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.03867933333333334
$ws.Range("H2").Value = 0.116038
$ws.Range("I2").Value = 0.2723398063265412
$ws.Range("J2").Value = 0.2723398063265412
$ws.Range("M2").Value = 13.35941066666667
$ws.Range("N2").Value = 40.078232
$ws.Range("O2").Value = 0.4925555025958562
$ws.Range("P2").Value = 0.4925555025958562
$ws.Range("Q2").Value = 0.5167330983128889
$ws.Range("R2").Value = 4.650597884816
$ws.Range("S2").Value = 0.1341424701820277
$ws.Range("T2").Value = 0.1341424701820276

# Row 3
$ws.Range("G3").Value = 0.03867933333333334
$ws.Range("H3").Value = 0.116038
$ws.Range("I3").Value = 0.2723398063265412
$ws.Range("J3").Value = 0.2723398063265412
$ws.Range("O3").Value = 0.03774352140193379
$ws.Range("P3").Value = 0.03774352140193379
$ws.Range("Q3").Value = 0.03959620114377778
$ws.Range("R3").Value = 0.356365810294
$ws.Range("S3").Value = 0.01027906330868431
$ws.Range("T3").Value = 0.01027906330868431

# Row 4
$ws.Range("G4").Value = 0.03867933333333334
$ws.Range("H4").Value = 0.116038
$ws.Range("I4").Value = 0.2723398063265412
$ws.Range("J4").Value = 0.2723398063265412
$ws.Range("M4").Value = 12.73953533333333
$ws.Range("N4").Value = 38.218606
$ws.Range("O4").Value = 0.4697009760022101
$ws.Range("P4").Value = 0.46970097600221
$ws.Range("Q4").Value = 0.4927567336697778
$ws.Range("R4").Value = 4.434810603028001
$ws.Range("S4").Value = 0.1279182728358293
$ws.Range("T4").Value = 0.1279182728358292

# Row 5
$ws.Range("I5").Value = 0.2603279211787514
$ws.Range("J5").Value = 0.2603279211787514
$ws.Range("M5").Value = 13.35941066666667
$ws.Range("N5").Value = 40.078232
$ws.Range("O5").Value = 0.4925555025958562
$ws.Range("P5").Value = 0.4925555025958562
$ws.Range("Q5").Value = 0.4939419437155556
$ws.Range("R5").Value = 4.44547749344
$ws.Range("S5").Value = 0.1282259500559343
$ws.Range("T5").Value = 0.1282259500559343

# Row 6
$ws.Range("I6").Value = 0.2603279211787514
$ws.Range("J6").Value = 0.2603279211787514
$ws.Range("O6").Value = 0.03774352140193379
$ws.Range("P6").Value = 0.03774352140193379
$ws.Range("S6").Value = 0.009825692464531136
$ws.Range("T6").Value = 0.009825692464531134

# Row 7
$ws.Range("I7").Value = 0.2603279211787514
$ws.Range("J7").Value = 0.2603279211787514
$ws.Range("M7").Value = 12.73953533333333
$ws.Range("N7").Value = 38.218606
$ws.Range("O7").Value = 0.4697009760022101
$ws.Range("P7").Value = 0.46970097600221
$ws.Range("Q7").Value = 0.4710230863911112
$ws.Range("R7").Value = 4.23920777752
$ws.Range("S7").Value = 0.1222762786582859
$ws.Range("T7").Value = 0.1222762786582859

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.06637333333333333
$ws.Range("H8").Value = 0.19912
$ws.Range("I8").Value = 0.4673322724947075
$ws.Range("J8").Value = 0.4673322724947075
$ws.Range("M8").Value = 13.35941066666667
$ws.Range("N8").Value = 40.078232
$ws.Range("O8").Value = 0.4925555025958562
$ws.Range("P8").Value = 0.4925555025958562
$ws.Range("Q8").Value = 0.8867086173155555
$ws.Range("R8").Value = 7.98037755584
$ws.Range("S8").Value = 0.2301870823578943
$ws.Range("T8").Value = 0.2301870823578943

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.06637333333333333
$ws.Range("H9").Value = 0.19912
$ws.Range("I9").Value = 0.4673322724947075
$ws.Range("J9").Value = 0.4673322724947075
$ws.Range("O9").Value = 0.03774352140193379
$ws.Range("P9").Value = 0.03774352140193379
$ws.Range("Q9").Value = 0.06794666895111111
$ws.Range("R9").Value = 0.61152002056
$ws.Range("S9").Value = 0.01763876562871835
$ws.Range("T9").Value = 0.01763876562871835

# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.06637333333333333
$ws.Range("H10").Value = 0.19912
$ws.Range("I10").Value = 0.4673322724947075
$ws.Range("J10").Value = 0.4673322724947075
$ws.Range("M10").Value = 12.73953533333333
$ws.Range("N10").Value = 38.218606
$ws.Range("O10").Value = 0.4697009760022101
$ws.Range("P10").Value = 0.46970097600221
$ws.Range("Q10").Value = 0.8455654251911111
$ws.Range("R10").Value = 7.61008882672
$ws.Range("S10").Value = 0.2195064245080949
$ws.Range("T10").Value = 0.2195064245080949

Write-Host "Applied TPM updates to rows 2-10"
